# Daily_Time_Series test cases.xlsx - apply commit changes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Row height tweaks on existing rows (content was reworded to be shorter)
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 75
$ws.Rows.Item(9).RowHeight = 75
$ws.Rows.Item(14).RowHeight = 120

# ---------------------------------------------------------------------------
# 2. Existing row 16 / 17 text corrections
# ---------------------------------------------------------------------------
$ws.Range("E16").Value = "1.Status code should be 200.`n2.Response should be in csv format"

$ws.Range("B17").Value = "Verify test case 4 and 8 with optional datatype parameter"
$ws.Range("D17").Value = "same as test case 4 and 8 appending datatype=csv to the get request"
$ws.Range("E17").Value = "1.Status code should be 200.`n2.Response should be in csv format"

# ---------------------------------------------------------------------------
# 3. New section header row 18 ("Meta Data")
# ---------------------------------------------------------------------------
$ws.Range("A18").Value = "Meta Data"
$ws.Range("A18").Font.Size = 11
$ws.Range("A18").Interior.Pattern = [Microsoft.Office.Interop.Excel.XlPattern]::xlPatternSolid
$ws.Range("A18").Interior.ThemeColor = [Microsoft.Office.Interop.Excel.XlThemeColor]::xlThemeColorAccent1
$ws.Range("A18").Interior.TintAndShade = 0.4

# ---------------------------------------------------------------------------
# 4. New test-case rows 19-24
# ---------------------------------------------------------------------------
$ws.Range("A19").Value = 14
$ws.Range("B19").Value = "Verify symbol in the response meta data should be IBM"
$ws.Range("C19").Value = "1.Should have end point uri`n2.Should have all required parameter details"
$ws.Range("D19").Value = "1. Send GET request to https://www.alphavantage.co/query?function=TIME_SERIES_DAILY&symbol=IBM&apikey=demo"
$ws.Range("E19").Value = "1.Status code should be 200.`n2.Response meta data symbol should show IBM."

$ws.Range("A20").Value = 15
$ws.Range("B20").Value = "Verify symbol in the response meta data should be TSCO.LON"
$ws.Range("C20").Value = "1.Should have end point uri`n2.Should have all required parameter details"
$ws.Range("D20").Value = "1. Send GET request to https://www.alphavantage.co/query?function=TIME_SERIES_DAILY&symbol=TSCO.LON&apikey=demo"
$ws.Range("E20").Value = "1.Status code should be 200.`n2.Response meta data symbol should show TSCO.LON."

$ws.Range("A21").Value = 16
$ws.Range("B21").Value = "Verify Output Size in the response meta data should be Full Size"
$ws.Range("C21").Value = "1.Should have end point uri`n2.Should have all required parameter details"
$ws.Range("D21").Value = "1.Send GET request to `nhttps://www.alphavantage.co/query?function=TIME_SERIES_DAILY&symbol=TSCO.LON&outputsize=full&apikey=demo"
$ws.Range("E21").Value = "1.Status code should be 200.`n2.Response meta data Output Size should show Full Size."

$ws.Range("A22").Value = 17
$ws.Range("B22").Value = "Verify Output Size in the response meta data should be Compact"
$ws.Range("C22").Value = "1.Should have end point uri`n2.Should have all required parameter details"
$ws.Range("D22").Value = "1.Send GET request to `nhttps://www.alphavantage.co/query?function=TIME_SERIES_DAILY&symbol=TSCO.LON&apikey=demo"
$ws.Range("E22").Value = "1.Status code should be 200.`n2.Response meta data Output Size should show Compact."

$ws.Range("A23").Value = 18
$ws.Range("B23").Value = "Verify Time Zone in the response meta data should be US/Eastern"
$ws.Range("C23").Value = "1.Should have end point uri`n2.Should have all required parameter details"
$ws.Range("D23").Value = "1.Send GET request to `nhttps://www.alphavantage.co/query?function=TIME_SERIES_DAILY&symbol=TSCO.LON&apikey=demo"
$ws.Range("E23").Value = "1.Status code should be 200.`n2.Response meta data Output Size should show US/Eastern."

$ws.Range("A24").Value = 19
$ws.Range("B24").Value = "Verify Information in the response meta data should be Daily Prices (open, high, low, close) and Volumes"
$ws.Range("C24").Value = "1.Should have end point uri`n2.Should have all required parameter details"
$ws.Range("D24").Value = "1.Send GET request to `nhttps://www.alphavantage.co/query?function=TIME_SERIES_DAILY&symbol=TSCO.LON&apikey=demo"
$ws.Range("E24").Value = "1.Status code should be 200.`n2.Response meta data Information should show Daily Prices (open, high, low, close) and Volumes."

# ---------------------------------------------------------------------------
# 5. New section header row 25 ("API Response")
# ---------------------------------------------------------------------------
$ws.Range("A25").Value = "API Response"
$ws.Range("A25:E25").Interior.Pattern = [Microsoft.Office.Interop.Excel.XlPattern]::xlPatternSolid
$ws.Range("A25:E25").Interior.ThemeColor = [Microsoft.Office.Interop.Excel.XlThemeColor]::xlThemeColorAccent1
$ws.Range("A25:E25").Interior.TintAndShade = 0.4

# ---------------------------------------------------------------------------
# 6. New test-case row 26
# ---------------------------------------------------------------------------
$ws.Range("A26").Value = 19
$ws.Range("B26").Value = "Response should contain Meta Data and Time Series (Daily)"
$ws.Range("C26").Value = "1.Should have end point uri`n2.Should have all required parameter details"
$ws.Range("D26").Value = "1.Send GET request to `nhttps://www.alphavantage.co/query?function=TIME_SERIES_DAILY&symbol=TSCO.LON&apikey=demo"
$ws.Range("E26").Value = "1.Status code should be 200.`n2.Response should contain Meta Data and Time Series (Daily)."

# ---------------------------------------------------------------------------
# 7. Formatting: wrap text + borders for the new rows
# ---------------------------------------------------------------------------
$ws.Range("B19:B24").WrapText = $true
$ws.Range("B26").WrapText = $true

$dataRange = $ws.Range("C19:E22")
$dataRange.Borders.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$dataRange.Borders.Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin
$dataRange.WrapText = $true

$borderRange = $ws.Range("C23:E23")
$borderRange.Borders.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$borderRange.Borders.Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin
$borderRange.Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom).LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlLineStyleNone
$borderRange.WrapText = $true

$ws.Range("C24:E24").Borders.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$ws.Range("C24:E24").Borders.Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin
$ws.Range("C24:E24").WrapText = $true

$ws.Range("C26:E26").Borders.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$ws.Range("C26:E26").Borders.Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin
$ws.Range("C26:E26").WrapText = $true

# Row heights for the new rows
$ws.Rows.Item(19).RowHeight = 75
$ws.Rows.Item(20).RowHeight = 75
$ws.Rows.Item(21).RowHeight = 75
$ws.Rows.Item(22).RowHeight = 75
$ws.Rows.Item(23).RowHeight = 75
$ws.Rows.Item(24).RowHeight = 75
$ws.Rows.Item(26).RowHeight = 75

# ---------------------------------------------------------------------------
# 8. Sheet view selection
# ---------------------------------------------------------------------------
$ws.Range("E33").Select()

# ---------------------------------------------------------------------------
# 9. Workbook MRU absolute path (x15ac:absPath) -- not exposed via the Excel
#    object model; left untouched.
# ---------------------------------------------------------------------------
